# Apply "Mutation fixed" commit:
#  - Operators sheet: mutate operator mj_random_jump -> mj_mutPolynomialBounded
#  - Parameters sheet: rework rows 4-8 into rows 4-7
#       Row4: Probability crossover  0.5 -> 0.9  (float)
#       Row5: Probability flip allele -> Probability mutation ; 0.5 -> 1 (float)
#       Row6: Jump size -> Crowding degree ; 10 -> 20 ; int -> float
#       Row7: Probability mutation -> indpb ; 0.5 -> 1/30 (float, 2-decimal display format)
#       Row8 (Crowding degree / 0.05) removed
#  - Selection on Parameters sheet moves to B7

$wb = $excel.ActiveWorkbook

# --- Operators sheet: change mutate operator ---
$wsOperators = $wb.Worksheets.Item("Operators")
$wsOperators.Range("C3").Value = "mj_mutPolynomialBounded"

# --- Parameters sheet: rebuild rows 4-7, drop old row 8 ---
$wsParams = $wb.Worksheets.Item("Parameters")

$wsParams.Range("A4").Value = "Probability crossover"
$wsParams.Range("B4").Value = 0.9
$wsParams.Range("C4").Value = "float"

$wsParams.Range("A5").Value = "Probability mutation"
$wsParams.Range("B5").Value = 1
$wsParams.Range("C5").Value = "float"

$wsParams.Range("A6").Value = "Crowding degree"
$wsParams.Range("B6").Value = 20
$wsParams.Range("C6").Value = "float"

$wsParams.Range("A7").Value = "indpb"
$wsParams.Range("B7").Value = 1/30
$wsParams.Range("B7").NumberFormat = "0.00"
$wsParams.Range("C7").Value = "float"

# Remove the now-unused 8th row entirely
$wsParams.Range("A8:C8").Delete()

# Match the new selection recorded in the workbook view
$wsParams.Range("B7").Select()
